$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.348.13"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.42"
$ws.Range("E3").Value = "  +1.09%  "

$ws.Range("E4").Value = "  +0.97%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.72"

$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4752"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3687"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07453"
$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8846"
$ws.Range("E10").Value = "  +1.92%  "

$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.96"
$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07320"
$ws.Range("E13").Value = "  +3.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.443"
$ws.Range("E14").Value = "  +1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.11"
$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.586"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.532.11"
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.66"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.095.67"
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.892"
$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.89"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.139"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.242"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.43"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09001"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7536"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.543"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.954"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05344"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.304"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.396"
$ws.Range("E42").Value = "  +4.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5314"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1659"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.473"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4912"
$ws.Range("E46").Value = "  +1.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.55"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.87"
$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.671"
$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06296"
$ws.Range("E51").Value = "  +0.12%  "
